$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.756.38"
$ws.Range("E2").Value = "  -5.67%  "

$ws.Range("D3").Value = "3.310.13"
$ws.Range("E3").Value = "  -6.91%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.17"
$ws.Range("E5").Value = "  -5.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.10"
$ws.Range("E6").Value = "  -8.65%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  -4.36%  "

$ws.Range("D9").Value = "3.301.93"
$ws.Range("E9").Value = "  -6.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.182"
$ws.Range("E10").Value = "  -11.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.580"
$ws.Range("E11").Value = "  -8.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.94"
$ws.Range("E12").Value = "  -11.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000261"
$ws.Range("E13").Value = "  -9.58%  "

$ws.Range("D14").Value = "3.847.21"
$ws.Range("E14").Value = "  -6.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.47"
$ws.Range("E15").Value = "  -8.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "594.68"
$ws.Range("E16").Value = "  -11.38%  "

$ws.Range("D17").Value = "65.744.00"
$ws.Range("E17").Value = "  -5.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.88"
$ws.Range("E18").Value = "  -3.40%  "

$ws.Range("E19").Value = "  -4.08%  "

$ws.Range("D20").Value = "3.311.78"
$ws.Range("E20").Value = "  -7.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.28"
$ws.Range("E21").Value = "  -10.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.897"
$ws.Range("E22").Value = "  -7.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.66"
$ws.Range("E23").Value = "  -7.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.04"
$ws.Range("E24").Value = "  -6.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.07"
$ws.Range("E25").Value = "  -6.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.96"
$ws.Range("E26").Value = "  -9.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.98"
$ws.Range("E27").Value = "  -0.55%  "

$ws.Range("E28").Value = "  -10.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.17"
$ws.Range("E29").Value = "  -10.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.58"
$ws.Range("E30").Value = "  -11.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.40"
$ws.Range("E31").Value = "  -9.42%  "

$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.79"
$ws.Range("E32").Value = "  -13.72%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.17"
$ws.Range("E33").Value = "  -9.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.90"
$ws.Range("E34").Value = "  -7.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.103"
$ws.Range("E35").Value = "  -7.48%  "

$ws.Range("D36").Value = "3.725.38"
$ws.Range("E36").Value = "  -1.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "57.39"
$ws.Range("E37").Value = "  -7.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "524.07"
$ws.Range("E38").Value = "  +4.35%  "

$ws.Range("E39").Value = "  -0.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.42"
$ws.Range("E40").Value = "  -8.51%  "

$ws.Range("D41").Value = "0.0₃0702"
$ws.Range("E41").Value = "  -14.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.63"
$ws.Range("E42").Value = "  -9.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.123"
$ws.Range("E43").Value = "  -8.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.336"
$ws.Range("E44").Value = "  -9.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.51"
$ws.Range("E45").Value = "  -9.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.27"
$ws.Range("E46").Value = "  -2.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0409"
$ws.Range("E47").Value = "  -9.19%  "

$ws.Range("E48").Value = "  +9.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.128"
$ws.Range("E49").Value = "  -6.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.58"
$ws.Range("E50").Value = "  -10.37%  "

$ws.Range("E51").Value = "  -0.07%  "
